$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Header date
Replace-Text "2026-02-15 Sunday" "2026-02-16 Monday"

# Table cells (order matters to avoid collisions between old/new values)
Replace-Text "417×2=" "822×4="
Replace-Text "446×3=" "226×2="
Replace-Text "187×5=" "253×6="
Replace-Text "288×2=" "901×3="
Replace-Text "361×3=" "964×7="
Replace-Text "438×3=" "319×5="
Replace-Text "600×2=" "343×5="
Replace-Text "965×4=" "256×2="
Replace-Text "555×8=" "379×5="
Replace-Text "832×3=" "726×7="
Replace-Text "965×6=" "640×6="
Replace-Text "467×5=" "384×7="
Replace-Text "128×9=" "601×3="
Replace-Text "302×5=" "348×5="
Replace-Text "885×7=" "965×6="
Replace-Text "200×9=" "284×4="
Replace-Text "322×6=" "823×4="
Replace-Text "104×2=" "479×3="
Replace-Text "777×7=" "524×3="
Replace-Text "216×6=" "971×7="
Replace-Text "757×4=" "261×5="
Replace-Text "651×9=" "431×5="
Replace-Text "733×5=" "846×2="
Replace-Text "565×4=" "366×3="
Replace-Text "858×8=" "150×8="
